$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (27 and 28) to the existing table, mirroring the
# structure of the other rows: Item, Multiplicador, Pontos, Batalhas,
# ValorApostado, ValorAcumulado, Resultado

$ws.Range("A27").Value = "BonusPower"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 760
$ws.Range("D27").Value = 96
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 200
$ws.Range("G27").Value = "win"

$ws.Range("A28").Value = "SkipBoss"
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 960
$ws.Range("D28").Value = 96
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 40
$ws.Range("G28").Value = "win"
